# Weekly forward look stats update 15.08.25
#
# 1. Bump the "as at" date in the intro line (A2).
# 2. Drop the now-past "11 Aug 2025 / Mortgage and landlord possession
#    statistics" row (old row 5) - the whole table shifts up by one.
# 3. Add the newly pre-announced "Ethnicity and the Criminal Justice
#    System 2024" publication for week commencing 24 Nov 2025, directly
#    under the other 24 Nov 2025 entry (HMPPS offender equalities report).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value2 = "This list contains a week-by-week view of  MoJ Official and National Statistics that have been pre-announced on the gov.uk release calendar as at 15 August 2025"

$ws.Rows(5).Delete()

$ws.Rows(28).Insert()
$ws.Range("A28").Value2 = "24 Nov 2025"
$ws.Range("B28").Value2 = "Ethnicity and the Criminal Justice System 2024"
$ws.Range("C28").Value2 = "27 November 2025"
$ws.Range("D28").Value2 = "provisional"
$ws.Range("E28").Value2 = 48
$ws.Range("F28").Value2 = "standard"
